# Update countries & provincias Spain
# - Refresh the "Datos actualizados" timestamp (13:22 -> 13:52)
# - Update Brasil's daily figures (row 14)
# - Insert Catar's new figures in its alphabetically-correct spot (row 41),
#   which pushes Chequia/Serbia/Filipinas/Australia/Ucrania down one row
#   each, and remove the old Catar row that used to sit after Ucrania.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 22 de Abril de 2020 a las 13:52"

# Brasil (row 14) daily figures
$ws.Range("B14").Value = 43592
$ws.Range("C14").Value = 513
$ws.Range("E14").Value = 16498
$ws.Range("G14").Value = 28
$ws.Range("H14").Value = 2769

# Insert a new row for Catar right before Chequia (row 41), shifting
# Chequia/Serbia/Filipinas/Australia/Ucrania down by one row each.
$ws.Rows.Item(41).Insert()

# Remove the now-duplicated old Catar row (previously right after Ucrania,
# now shifted down to row 47).
$ws.Rows.Item(47).Delete()

# Populate the new Catar row with its updated figures.
$ws.Range("A41").Value = "Catar"
$ws.Range("B41").Value = 7141
$ws.Range("C41").Value = 608
$ws.Range("D41").Value = 689
$ws.Range("E41").Value = 6442
$ws.Range("F41").Value = 37
$ws.Range("G41").Value = 1
$ws.Range("H41").Value = 10
